# Add new scheduler responses (Julian, Maddy, Emily, Ada) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting for the new rows by copying the style from the existing
# data row (row 2) down into rows 3-6.
$ws.Range("A2:K2").Copy()
$ws.Range("A3:K6").PasteSpecial(-4122)

$ws.Range("A3").Value = 44752.82988744213
$ws.Range("B3").Value = "Julian"
$ws.Range("C3").Value = "8 AM - 9 AM, 11 AM - 12 PM, 1 PM - 2 PM, 7 PM - 8 PM"
$ws.Range("D3").Value = "10 AM - 11 AM, 1 PM - 2 PM, 6 PM - 7 PM, 7 PM - 8 PM, 8 PM - 9 PM, 10 PM - 11 PM"
$ws.Range("E3").Value = "10 AM - 11 AM, 11 AM - 12 PM, 12 PM - 1 PM, 1 PM - 2 PM, 3 PM - 4 PM, 5 PM - 6 PM, 6 PM - 7 PM, 7 PM - 8 PM, 8 PM - 9 PM, 9 PM - 10 PM, 10 PM - 11 PM, 11 PM - 12 AM"
$ws.Range("F3").Value = "9 AM - 10 AM, 10 AM - 11 AM, 12 PM - 1 PM, 1 PM - 2 PM, 5 PM - 6 PM, 7 PM - 8 PM"
$ws.Range("G3").Value = "9 AM - 10 AM, 10 AM - 11 AM, 11 AM - 12 PM, 12 PM - 1 PM, 1 PM - 2 PM, 3 PM - 4 PM, 6 PM - 7 PM, 7 PM - 8 PM, 8 PM - 9 PM, 9 PM - 10 PM, 10 PM - 11 PM"
$ws.Range("H3").Value = "8 AM - 9 AM, 10 AM - 11 AM, 11 AM - 12 PM, 12 PM - 1 PM, 6 PM - 7 PM, 7 PM - 8 PM, 10 PM - 11 PM, 11 PM - 12 AM"
$ws.Range("I3").Value = "10 AM - 11 AM, 12 PM - 1 PM, 5 PM - 6 PM, 6 PM - 7 PM, 7 PM - 8 PM, 9 PM - 10 PM, 11 PM - 12 AM"
$ws.Range("J3").Value = "Shop, Fridges, Bathrooms, Dining Room Clean, Trash & Recycling, Garden Helper"
$ws.Range("K3").Value = "Fridges, Compost, Dining Room Clean, First Floor Commons Clean"
$ws.Range("A4").Value = 44752.83030295139
$ws.Range("B4").Value = "Maddy"
$ws.Range("C4").Value = "8 AM - 9 AM, 9 AM - 10 AM, 10 AM - 11 AM, 11 AM - 12 PM, 12 PM - 1 PM, 7 PM - 8 PM, 10 PM - 11 PM, 11 PM - 12 AM"
$ws.Range("D4").Value = "11 AM - 12 PM, 3 PM - 4 PM, 5 PM - 6 PM, 6 PM - 7 PM, 8 PM - 9 PM, 9 PM - 10 PM, 10 PM - 11 PM, 11 PM - 12 AM"
$ws.Range("E4").Value = "9 AM - 10 AM, 10 AM - 11 AM, 12 PM - 1 PM, 1 PM - 2 PM, 3 PM - 4 PM, 8 PM - 9 PM, 9 PM - 10 PM, 10 PM - 11 PM"
$ws.Range("F4").Value = "9 AM - 10 AM, 10 AM - 11 AM, 12 PM - 1 PM, 1 PM - 2 PM, 5 PM - 6 PM, 8 PM - 9 PM, 10 PM - 11 PM"
$ws.Range("G4").Value = "9 AM - 10 AM, 12 PM - 1 PM, 1 PM - 2 PM, 6 PM - 7 PM, 9 PM - 10 PM, 10 PM - 11 PM"
$ws.Range("H4").Value = "10 AM - 11 AM, 11 AM - 12 PM, 12 PM - 1 PM, 1 PM - 2 PM, 3 PM - 4 PM, 5 PM - 6 PM, 6 PM - 7 PM, 8 PM - 9 PM, 10 PM - 11 PM"
$ws.Range("I4").Value = "9 AM - 10 AM, 11 AM - 12 PM, 12 PM - 1 PM, 5 PM - 6 PM, 6 PM - 7 PM, 7 PM - 8 PM"
$ws.Range("J4").Value = "Fast Cook, Lunch Clean, Kitchen Deep Clean, Compost, First Floor Commons Clean, Porch Yard Clean, Laundry Room and Rags, Garden Helper"
$ws.Range("K4").Value = "Dinner Cook, Fast Cook, Fridges, Compost, Bathrooms, Porch Yard Clean"
$ws.Range("A5").Value = 44752.830715648146
$ws.Range("B5").Value = "Emily"
$ws.Range("C5").Value = "8 AM - 9 AM, 9 AM - 10 AM, 10 AM - 11 AM, 12 PM - 1 PM, 1 PM - 2 PM, 3 PM - 4 PM, 7 PM - 8 PM, 9 PM - 10 PM, 10 PM - 11 PM, 11 PM - 12 AM"
$ws.Range("D5").Value = "8 AM - 9 AM, 9 AM - 10 AM, 10 AM - 11 AM, 11 AM - 12 PM, 12 PM - 1 PM, 1 PM - 2 PM, 3 PM - 4 PM, 5 PM - 6 PM, 6 PM - 7 PM, 7 PM - 8 PM, 8 PM - 9 PM, 9 PM - 10 PM, 10 PM - 11 PM, 11 PM - 12 AM"
$ws.Range("E5").Value = "1 PM - 2 PM, 5 PM - 6 PM, 6 PM - 7 PM, 8 PM - 9 PM, 10 PM - 11 PM"
$ws.Range("F5").Value = "1 PM - 2 PM, 3 PM - 4 PM, 5 PM - 6 PM, 6 PM - 7 PM, 7 PM - 8 PM, 8 PM - 9 PM, 10 PM - 11 PM"
$ws.Range("G5").Value = "9 AM - 10 AM, 12 PM - 1 PM, 5 PM - 6 PM, 7 PM - 8 PM, 8 PM - 9 PM, 9 PM - 10 PM, 11 PM - 12 AM"
$ws.Range("H5").Value = "8 AM - 9 AM, 9 AM - 10 AM, 10 AM - 11 AM, 3 PM - 4 PM, 5 PM - 6 PM, 7 PM - 8 PM, 8 PM - 9 PM, 9 PM - 10 PM, 10 PM - 11 PM"
$ws.Range("I5").Value = "9 AM - 10 AM, 10 AM - 11 AM, 12 PM - 1 PM, 8 PM - 9 PM"
$ws.Range("J5").Value = "Dinner Cook, Fast Cook, Dinner Clean, Fridges, Kitchen Deep Clean, Bathrooms, Dining Room Clean"
$ws.Range("K5").Value = "Lunch Clean, Kitchen Deep Clean, Dining Room Clean, Pool Clean"
$ws.Range("A6").Value = 44752.83098380787
$ws.Range("B6").Value = "Ada"
$ws.Range("C6").Value = "9 AM - 10 AM, 11 AM - 12 PM, 12 PM - 1 PM, 6 PM - 7 PM"
$ws.Range("D6").Value = "11 AM - 12 PM"
$ws.Range("E6").Value = "9 AM - 10 AM, 12 PM - 1 PM, 3 PM - 4 PM, 6 PM - 7 PM, 7 PM - 8 PM, 8 PM - 9 PM, 9 PM - 10 PM, 10 PM - 11 PM, 11 PM - 12 AM"
$ws.Range("F6").Value = "8 AM - 9 AM, 10 AM - 11 AM, 11 AM - 12 PM, 1 PM - 2 PM, 6 PM - 7 PM, 9 PM - 10 PM, 10 PM - 11 PM"
$ws.Range("G6").Value = "12 PM - 1 PM, 1 PM - 2 PM, 5 PM - 6 PM, 8 PM - 9 PM"
$ws.Range("H6").Value = "11 AM - 12 PM, 12 PM - 1 PM, 6 PM - 7 PM, 7 PM - 8 PM, 8 PM - 9 PM, 10 PM - 11 PM"
$ws.Range("I6").Value = "3 PM - 4 PM, 6 PM - 7 PM, 8 PM - 9 PM"
$ws.Range("J6").Value = "Fast Cook, Dinner Clean, Kitchen Deep Clean, First Floor Commons Clean, Trash & Recycling, Garden Helper"
$ws.Range("K6").Value = "Shop, Kitchen Deep Clean, Bathrooms"
